$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Range("H2").Value = 256.8
$ws.Range("I2").Value = 154.5
$ws.Range("K2").Value = 154.5
$ws.Range("M2").Value = -41.5
# Row 92
$ws.Range("H92").Value = 2980.6365
$ws.Range("I92").Value = 3387.4443
$ws.Range("J92").Value = 1150
$ws.Range("K92").Value = 3387.4443
$ws.Range("L92").Value = 1150
$ws.Range("M92").Value = -2139.4443
$ws.Range("N92").Value = -3646
# Row 106
$ws.Range("H106").Value = 13061.8
$ws.Range("I106").Value = 3827.5
$ws.Range("J106").Value = 49999
$ws.Range("K106").Value = 3827.5
$ws.Range("L106").Value = 49999
$ws.Range("M106").Value = -3196.5
$ws.Range("N106").Value = -51261
# Row 137
$ws.Range("H137").Value = 9924.628000000001
$ws.Range("J137").Value = 15469.583
$ws.Range("L137").Value = 46408.749
$ws.Range("N137").Value = -51508.749

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 10855.354
$ws.Range("I32").Value = 3942.4746
$ws.Range("K32").Value = 3942.4746
$ws.Range("M32").Value = -3655.4746
# Row 61
$ws.Range("H61").Value = 5994.6665
$ws.Range("I61").Value = 1648.7858
$ws.Range("K61").Value = 1648.7858
$ws.Range("M61").Value = -1436.7858
# Row 68
$ws.Range("H68").Value = 60033
$ws.Range("J68").Value = 60033
$ws.Range("L68").Value = 60033
$ws.Range("N68").Value = -61655
# Row 71
$ws.Range("H71").Value = 60033
$ws.Range("J71").Value = 60033
$ws.Range("L71").Value = 180099
$ws.Range("N71").Value = -188211
# Row 74
$ws.Range("H74").Value = 16154.552
$ws.Range("I74").Value = 1978.2307
$ws.Range("K74").Value = 1978.2307
$ws.Range("M74").Value = -1104.2307
# Row 77
$ws.Range("H77").Value = 16154.552
$ws.Range("I77").Value = 1978.2307
$ws.Range("K77").Value = 9891.1535
$ws.Range("M77").Value = -5523.1535
# Row 110
$ws.Range("H110").Value = 10539.842
$ws.Range("I110").Value = 6337.4
$ws.Range("K110").Value = 6337.4
$ws.Range("M110").Value = -4292.4
# Row 132
$ws.Range("H132").Value = 1858843.1
$ws.Range("I132").Value = 2549.0698
$ws.Range("J132").Value = 9115266
$ws.Range("K132").Value = 7647.209400000001
$ws.Range("L132").Value = 27345798
$ws.Range("M132").Value = -5117.209400000001
$ws.Range("N132").Value = -27350858
# Row 136
$ws.Range("H136").Value = 5994.6665
$ws.Range("I136").Value = 1648.7858
$ws.Range("K136").Value = 4946.357400000001
$ws.Range("M136").Value = -2396.357400000001

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 22263.926
$ws.Range("I20").Value = 12069.606
$ws.Range("K20").Value = 12069.606
$ws.Range("M20").Value = -11822.606
# Row 94
$ws.Range("H94").Value = 1775.6
$ws.Range("J94").Value = 4816.875
$ws.Range("L94").Value = 4816.875
$ws.Range("N94").Value = -5718.875
# Row 107
$ws.Range("H107").Value = 1856.4166
$ws.Range("I107").Value = 1587.7
$ws.Range("J107").Value = 3200
$ws.Range("K107").Value = 1587.7
$ws.Range("L107").Value = 3200
$ws.Range("M107").Value = 332.3
$ws.Range("N107").Value = -7040
# Row 124
$ws.Range("H124").Value = 39254.5
$ws.Range("J124").Value = 39254.5
$ws.Range("L124").Value = 39254.5
$ws.Range("N124").Value = -49074.5
# Row 132
$ws.Range("H132").Value = 93446.3
$ws.Range("J132").Value = 93446.3
$ws.Range("L132").Value = 93446.3
$ws.Range("N132").Value = -103566.3
# Row 134
$ws.Range("H134").Value = 14288.777
$ws.Range("I134").Value = 7401.0527
$ws.Range("K134").Value = 22203.1581
$ws.Range("M134").Value = -19668.1581

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 10975.551
$ws.Range("I58").Value = 5594.0625
$ws.Range("J58").Value = 13584.758
$ws.Range("K58").Value = 5594.0625
$ws.Range("L58").Value = 13584.758
$ws.Range("M58").Value = -5391.0625
$ws.Range("N58").Value = -13990.758
# Row 62
$ws.Range("H62").Value = 5022.6
$ws.Range("I62").Value = 4855
$ws.Range("J62").Value = 5190.2
$ws.Range("K62").Value = 4855
$ws.Range("L62").Value = 5190.2
$ws.Range("M62").Value = -4231
$ws.Range("N62").Value = -6438.2
# Row 65
$ws.Range("H65").Value = 5022.6
$ws.Range("I65").Value = 4855
$ws.Range("J65").Value = 5190.2
$ws.Range("K65").Value = 24275
$ws.Range("L65").Value = 25951
$ws.Range("M65").Value = -21155
$ws.Range("N65").Value = -32191
# Row 134
$ws.Range("H134").Value = 23261244
$ws.Range("I134").Value = 1621.8182
$ws.Range("K134").Value = 4865.4546
$ws.Range("M134").Value = -2330.4546
# Row 136
$ws.Range("H136").Value = 10975.551
$ws.Range("I136").Value = 5594.0625
$ws.Range("J136").Value = 13584.758
$ws.Range("K136").Value = 16782.1875
$ws.Range("L136").Value = 40754.274
$ws.Range("M136").Value = -14232.1875
$ws.Range("N136").Value = -45854.274

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 1436096
$ws.Range("I5").Value = 1008.9091
$ws.Range("K5").Value = 3026.7273
$ws.Range("M5").Value = -2914.7273
# Row 12
$ws.Range("H12").Value = 28.272728
$ws.Range("J12").Value = 38.5
$ws.Range("L12").Value = 115.5
$ws.Range("N12").Value = -461.5
# Row 42
$ws.Range("H42").Value = 8000
$ws.Range("J42").Value = 8000
$ws.Range("L42").Value = 24000
$ws.Range("N42").Value = -25068
# Row 104
$ws.Range("H104").Value = 1373153.2
$ws.Range("I104").Value = 1587.5
$ws.Range("K104").Value = 4762.5
$ws.Range("M104").Value = -2141.5
# Row 107
$ws.Range("H107").Value = 2083982.4
$ws.Range("J107").Value = 3906966
$ws.Range("L107").Value = 11720898
$ws.Range("N107").Value = -11724738
# Row 135
$ws.Range("H135").Value = 1436096
$ws.Range("I135").Value = 1008.9091
$ws.Range("K135").Value = 9080.1819
$ws.Range("M135").Value = -6545.1819

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 19
$ws.Range("H19").Value = 2500
$ws.Range("J19").Value = 2500
$ws.Range("L19").Value = 2500
$ws.Range("N19").Value = -3076
# Row 70
$ws.Range("H70").Value = 37713.715
$ws.Range("I70").Value = 33199.4
$ws.Range("K70").Value = 33199.4
$ws.Range("M70").Value = -32929.4
# Row 73
$ws.Range("H73").Value = 37713.715
$ws.Range("I73").Value = 33199.4
$ws.Range("K73").Value = 33199.4
$ws.Range("M73").Value = -32263.4
# Row 80
$ws.Range("H80").Value = 18038.809
$ws.Range("I80").Value = 12958.7
$ws.Range("J80").Value = 22657.092
$ws.Range("K80").Value = 12958.7
$ws.Range("L80").Value = 22657.092
$ws.Range("M80").Value = -11960.7
$ws.Range("N80").Value = -24653.092
# Row 83
$ws.Range("H83").Value = 18038.809
$ws.Range("I83").Value = 12958.7
$ws.Range("J83").Value = 22657.092
$ws.Range("K83").Value = 64793.5
$ws.Range("L83").Value = 113285.46
$ws.Range("M83").Value = -59801.5
$ws.Range("N83").Value = -123269.46
# Row 107
$ws.Range("H107").Value = 650.3333
$ws.Range("J107").Value = 1094.2858
$ws.Range("L107").Value = 1094.2858
$ws.Range("N107").Value = -4934.2858
# Row 113
$ws.Range("H113").Value = 51844
$ws.Range("I113").Value = 66858.664
$ws.Range("J113").Value = 6800
$ws.Range("K113").Value = 66858.664
$ws.Range("L113").Value = 6800
$ws.Range("M113").Value = -64688.664
$ws.Range("N113").Value = -11140
# Row 123
$ws.Range("H123").Value = 50326
$ws.Range("J123").Value = 50326
$ws.Range("L123").Value = 50326
$ws.Range("N123").Value = -55226
# Row 132
$ws.Range("H132").Value = 3631.8909
$ws.Range("I132").Value = 1481.4546
$ws.Range("J132").Value = 12233.637
$ws.Range("K132").Value = 4444.3638
$ws.Range("L132").Value = 36700.911
$ws.Range("M132").Value = -1914.3638
$ws.Range("N132").Value = -41760.911

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 61
$ws.Range("H61").Value = 3481.1482
$ws.Range("I61").Value = 1858.25
$ws.Range("K61").Value = 1858.25
$ws.Range("M61").Value = -1656.25
# Row 68
$ws.Range("H68").Value = 9289.066000000001
$ws.Range("J68").Value = 9793.857
$ws.Range("L68").Value = 9793.857
$ws.Range("N68").Value = -11291.857
# Row 71
$ws.Range("H71").Value = 9289.066000000001
$ws.Range("J71").Value = 9793.857
$ws.Range("L71").Value = 48969.285
$ws.Range("N71").Value = -56457.285
# Row 87
$ws.Range("H87").Value = 100189
$ws.Range("J87").Value = 100189
$ws.Range("L87").Value = 100189
$ws.Range("N87").Value = -102435
# Row 90
$ws.Range("H90").Value = 100189
$ws.Range("J90").Value = 100189
$ws.Range("L90").Value = 300567
$ws.Range("N90").Value = -311799
# Row 113
$ws.Range("H113").Value = 3481.1482
$ws.Range("I113").Value = 1858.25
$ws.Range("K113").Value = 1858.25
$ws.Range("M113").Value = 311.75
# Row 125
$ws.Range("H125").Value = 129800
$ws.Range("J125").Value = 129800
$ws.Range("L125").Value = 129800
$ws.Range("N125").Value = -139640
# Row 136
$ws.Range("H136").Value = 18580.305
$ws.Range("I136").Value = 20403.584
$ws.Range("J136").Value = 16591.273
$ws.Range("K136").Value = 61210.75199999999
$ws.Range("L136").Value = 49773.819
$ws.Range("M136").Value = -58660.75199999999
$ws.Range("N136").Value = -54873.819

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 74
$ws.Range("H74").Value = 8525
$ws.Range("J74").Value = 8525
$ws.Range("L74").Value = 8525
$ws.Range("N74").Value = -10397
# Row 77
$ws.Range("H77").Value = 8525
$ws.Range("J77").Value = 8525
$ws.Range("L77").Value = 25575
$ws.Range("N77").Value = -34935
# Row 136
$ws.Range("H136").Value = 7503.39
$ws.Range("I136").Value = 1676.8966
$ws.Range("K136").Value = 5030.6898
$ws.Range("M136").Value = -2480.6898
